# Commit: "commit on 22 july 2019"
#
# Update the RegSheet test-fixture row 1 to hold new registration values and
# add a boolean flag row below it:
#   A1: "writ"     -> "helloworld1"
#   B1: "hi"       -> "test1"
#   C1: "newValue" -> (removed)
#   A2: (new)      -> FALSE

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "helloworld1"
$ws.Range("B1").Value = "test1"
$ws.Range("C1").ClearContents()
$ws.Range("A2").Value = $false
